$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B6's text value (shared string "foo bar nothing" -> "foo bar")
$ws.Range("B6").Value = "foo bar"

# Move the active selection from D7 to B7
[void]$ws.Range("B7").Select()
